# Add tutorial links to schedule
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix ellipsis prefix on "continued: Putting it all together" rows
$ws.Range("D41").Value = "… continued: Putting it all together"
$ws.Range("D43").Value = "… continued: Putting it all together"
$ws.Range("D45").Value = "… continued: Putting it all together"

# Add tutorial page links (column H, link_lab) alongside each lecture slide link
$ws.Range("H6").Value  = "pages/git.html"
$ws.Range("H8").Value  = "pages/git.html"

$ws.Range("H10").Value = "pages/conda.html"

$ws.Range("H12").Value = "pages/snakemake.html"
$ws.Range("H14").Value = "pages/snakemake.html"
$ws.Range("H16").Value = "pages/snakemake.html"

$ws.Range("H18").Value = "pages/nextflow.html"
$ws.Range("H21").Value = "pages/nextflow.html"
$ws.Range("H23").Value = "pages/nextflow.html"

$ws.Range("H25").Value = "pages/quarto.html"
$ws.Range("H27").Value = "pages/quarto.html"

$ws.Range("H28").Value = "pages/jupyter.html"
$ws.Range("H30").Value = "pages/jupyter.html"

$ws.Range("H31").Value = "pages/containers.html"
$ws.Range("H33").Value = "pages/containers.html"
$ws.Range("H35").Value = "pages/containers.html"
$ws.Range("H37").Value = "pages/containers.html"

# Widen column H to fit the new links
$ws.Columns("H").ColumnWidth = 19.7

# Update the selected cell to reflect where the author ended up
$ws.Range("G35").Select()
